$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The second "PRAZO DE VIGENCIA DO PRESENTE CONTRATO" heading is
#    retitled to "DO FORO" and the (hidden) "_GoBack" bookmark that
#    used to sit near the signature block is moved to the start of
#    this paragraph (Word keeps "_GoBack" unique, so re-adding it
#    here automatically drops the old one).
# ------------------------------------------------------------------
$headingCount = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "PRAZO DE VIG*NCIA DO PRESENTE CONTRATO*") {
        $headingCount = $headingCount + 1
        if ($headingCount -eq 2) {
            $startRng = $p.Range.Duplicate
            $startRng.Collapse(1)
            [void]$d.Bookmarks.Add("_GoBack", $startRng)

            [void]$p.Range.Find.Execute("PRAZO DE VIGÊNCIA DO PRESENTE CONTRATO", $true, $false, $false, $false, $false, $true, 1, $false, "DO FORO", 2)
        }
    }
}

# ------------------------------------------------------------------
# 2) Drop the centered alignment of the blank paragraph that sits
#    right after "Para dirimir quaisquer controversias..." (it is
#    the second blank line following that sentence).
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Para dirimir quaisquer controv*rsias oriundas do contrato*") {
        $blank1 = $p.Next()
        $blank2 = $blank1.Next()
        $blank2.Alignment = 0
    }
}
